$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap paired rows: the re-scrape reordered several adjacent match rows. ---
# --- Column A (row index) and the league/date columns (C/D/E) are untouched; ---
# --- everything else (B, F:AC) is exchanged between the two rows.           ---
$swapCols = @("B","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Row($ws, $r1, $r2, $cols) {
    foreach ($col in $cols) {
        $a = $ws.Range("$col$r1").Value2
        $b = $ws.Range("$col$r2").Value2
        $ws.Range("$col$r1").Value2 = $b
        $ws.Range("$col$r2").Value2 = $a
    }
}

Swap-Row $ws 9 10 $swapCols
Swap-Row $ws 29 30 $swapCols
Swap-Row $ws 76 77 $swapCols
Swap-Row $ws 87 88 $swapCols
Swap-Row $ws 111 112 $swapCols
Swap-Row $ws 122 123 $swapCols

# --- Append 3 new match rows (158-160). Copy A157:AC157 formatting first (one
#     Copy per PasteSpecial, CutCopyMode cleared only once at the end so the
#     per-column style mapping -- bold/border on col A, date format on col E --
#     is preserved), then overwrite the values.
$ws.Range("A157:AC157").Copy()
$ws.Range("A158:AC158").PasteSpecial(-4122)
$ws.Range("A157:AC157").Copy()
$ws.Range("A159:AC159").PasteSpecial(-4122)
$ws.Range("A157:AC157").Copy()
$ws.Range("A160:AC160").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 158
$ws.Range("A158").Value2 = 156
$ws.Range("B158").Value2 = 7952751
$ws.Range("C158").Value2 = "Bosnia Herzegovina Premier Liga"
$ws.Range("D158").Value2 = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E158").Value2 = 45401.45833333334
$ws.Range("F158").Value2 = "GOSK Gabela"
$ws.Range("G158").Value2 = "FK Tuzla City"
$ws.Range("H158").Value2 = 1
$ws.Range("I158").Value2 = 0
$ws.Range("J158").Value2 = "H"
$ws.Range("K158").Value2 = 1.666
$ws.Range("L158").Value2 = 4
$ws.Range("M158").Value2 = 3.75
$ws.Range("N158").Value2 = 2.375
$ws.Range("O158").Value2 = 3.8
$ws.Range("P158").Value2 = 2.375
$ws.Range("Q158").Value2 = 0
$ws.Range("R158").Value2 = 1.9
$ws.Range("S158").Value2 = 1.9
$ws.Range("T158").Value2 = 2.75
$ws.Range("U158").Value2 = 1.95
$ws.Range("V158").Value2 = 1.85
$ws.Range("W158").Value2 = 1.375
$ws.Range("X158").Value2 = -1
$ws.Range("Y158").Value2 = -1
$ws.Range("Z158").Value2 = 0.8999999999999999
$ws.Range("AA158").Value2 = -1
$ws.Range("AB158").Value2 = -1
$ws.Range("AC158").Value2 = 0.8500000000000001

# Row 159
$ws.Range("A159").Value2 = 157
$ws.Range("B159").Value2 = 7952754
$ws.Range("C159").Value2 = "Bosnia Herzegovina Premier Liga"
$ws.Range("D159").Value2 = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E159").Value2 = 45401.54166666666
$ws.Range("F159").Value2 = "NK Posusje"
$ws.Range("G159").Value2 = "Zeljeznicar"
$ws.Range("H159").Value2 = 1
$ws.Range("I159").Value2 = 0
$ws.Range("J159").Value2 = "H"
$ws.Range("K159").Value2 = 1.727
$ws.Range("L159").Value2 = 3.5
$ws.Range("M159").Value2 = 4
$ws.Range("N159").Value2 = 5
$ws.Range("O159").Value2 = 4
$ws.Range("P159").Value2 = 1.55
$ws.Range("Q159").Value2 = 1
$ws.Range("R159").Value2 = 1.8
$ws.Range("S159").Value2 = 2
$ws.Range("T159").Value2 = 2.25
$ws.Range("U159").Value2 = 1.8
$ws.Range("V159").Value2 = 2
$ws.Range("W159").Value2 = 4
$ws.Range("X159").Value2 = -1
$ws.Range("Y159").Value2 = -1
$ws.Range("Z159").Value2 = 0.8
$ws.Range("AA159").Value2 = -1
$ws.Range("AB159").Value2 = -1
$ws.Range("AC159").Value2 = 1

# Row 160
$ws.Range("A160").Value2 = 158
$ws.Range("B160").Value2 = 7952750
$ws.Range("C160").Value2 = "Bosnia Herzegovina Premier Liga"
$ws.Range("D160").Value2 = "Bosnia  Herzegovina Premier Liga"
$ws.Range("E160").Value2 = 45401.625
$ws.Range("F160").Value2 = "Velez Mostar"
$ws.Range("G160").Value2 = "Zvijezda 09"
$ws.Range("H160").Value2 = 3
$ws.Range("I160").Value2 = 1
$ws.Range("J160").Value2 = "H"
$ws.Range("K160").Value2 = 1.25
$ws.Range("L160").Value2 = 5.75
$ws.Range("M160").Value2 = 7
$ws.Range("N160").Value2 = 1.166
$ws.Range("O160").Value2 = 6
$ws.Range("P160").Value2 = 13
$ws.Range("Q160").Value2 = -2
$ws.Range("R160").Value2 = 1.875
$ws.Range("S160").Value2 = 1.925
$ws.Range("T160").Value2 = 3.25
$ws.Range("U160").Value2 = 1.95
$ws.Range("V160").Value2 = 1.85
$ws.Range("W160").Value2 = 0.1659999999999999
$ws.Range("X160").Value2 = -1
$ws.Range("Y160").Value2 = -1
$ws.Range("Z160").Value2 = 0
$ws.Range("AA160").Value2 = -0.0
$ws.Range("AB160").Value2 = 0.95
$ws.Range("AC160").Value2 = -1

Write-Host "edit applied"
